$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D/E) in this sheet store plain text cells
# (numbers and percentages are authored as literal strings, not real
# Number/Percentage cells). Pre-marking each touched D/E cell as Text
# ("@") before assigning its new value keeps Excel from silently
# re-typing strings like "246.94" or "1.03%" into numeric values.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

# Apply the updated cell values (coin name/link moved rows plus
# refreshed price/volume figures), row by row in document order.

# Row 2
$ws.Range("D2").Value = "246.94"
$ws.Range("E2").Value = "1.03%"

# Row 3
$ws.Range("D3").Value = "29.91"
$ws.Range("E3").Value = "10.34%"

# Row 4
$ws.Range("D4").Value = "5.173"
$ws.Range("E4").Value = "0.67%"

# Row 5
$ws.Range("D5").Value = "0.05719"
$ws.Range("E5").Value = "1.13%"

# Row 6
$ws.Range("D6").Value = "6.613"
$ws.Range("E6").Value = "2.18%"

# Row 7
$ws.Range("D7").Value = "3.053"
$ws.Range("E7").Value = "1.66%"

# Row 8
$ws.Range("D8").Value = "0.8597"
$ws.Range("E8").Value = "4.70%"

# Row 9
$ws.Range("D9").Value = "0.8695"
$ws.Range("E9").Value = "3.50%"

# Row 10
$ws.Range("D10").Value = "0.1361"
$ws.Range("E10").Value = "2.35%"

# Row 11
$ws.Range("D11").Value = "0.07068"
$ws.Range("E11").Value = "1.99%"

# Row 12
$ws.Range("D12").Value = "0.02866"
$ws.Range("E12").Value = "-3.76%"

# Row 13
$ws.Range("D13").Value = "0.09398"
$ws.Range("E13").Value = "0.05%"

# Row 14
$ws.Range("D14").Value = "0.001524"
$ws.Range("E14").Value = "0.12%"

# Row 15
$ws.Range("D15").Value = "0.04151"
$ws.Range("E15").Value = "-2.68%"

# Row 16
$ws.Range("D16").Value = "0.0005980"
$ws.Range("E16").Value = "-0.02%"

# Row 17
$ws.Range("D17").Value = "0.006136"
$ws.Range("E17").Value = "-0.09%"

# Row 18
$ws.Range("B18").Value = "UpBots"
$ws.Range("C18").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D18").Value = "0.007491"
$ws.Range("E18").Value = "10,219.03%"

# Row 19
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "3.478"
$ws.Range("E19").Value = "-1.06%"

# Row 20
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "2.171"
$ws.Range("E20").Value = "-2.57%"

# Row 21
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "0.3145"
$ws.Range("E21").Value = "1.03%"

# Row 22
$ws.Range("B22").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C22").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D22").Value = "0.03252"
$ws.Range("E22").Value = "3.52%"

# Row 23
$ws.Range("B23").Value = "ProBitToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D23").Value = "0.1300"
$ws.Range("E23").Value = "3.53%"

# Row 24
$ws.Range("B24").Value = "MCDex"
$ws.Range("C24").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D24").Value = "3.144"
$ws.Range("E24").Value = "-12.17%"

# Row 25
$ws.Range("B25").Value = "ZBToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D25").Value = "0.1380"
$ws.Range("E25").Value = "0.44%"

# Row 26
$ws.Range("B26").Value = "HotbitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").Value = "0.005092"
$ws.Range("E26").Value = "14.01%"

# Row 27
$ws.Range("B27").Value = "BitKan"
$ws.Range("C27").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D27").Value = "0.001220"
$ws.Range("E27").Value = "-0.15%"

# Row 28
$ws.Range("B28").Value = "NitroEx"
$ws.Range("C28").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D28").Value = "0.0001210"
$ws.Range("E28").Value = "23.44%"

# Row 40
$ws.Range("D40").Value = "0.03775"
$ws.Range("E40").Value = "3.36%"

# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1071"
$ws.Range("E41").Value = "1.57%"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002600"
$ws.Range("E42").Value = "-0.40%"

# Row 43
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003497"
$ws.Range("E43").Value = "-42.20%"

# Row 44
$ws.Range("D44").Value = "0.009634"
$ws.Range("E44").Value = "7.19%"

# Row 45
$ws.Range("D45").Value = "0.00005091"
$ws.Range("E45").Value = "-5.21%"

# Row 46
$ws.Range("E46").Value = "-0.02%"

# Row 47
$ws.Range("D47").Value = "0.07510"
$ws.Range("E47").Value = "-25.66%"

# Row 48
$ws.Range("D48").Value = "0.002758"
$ws.Range("E48").Value = "3.92%"

# Row 49
$ws.Range("E49").Value = "-0.02%"

# Row 50
$ws.Range("E50").Value = "-0.02%"
